# Auto-generated edit script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.659.19"
$ws.Range("E2").Value = "  -4.32%  "
$ws.Range("D3").Value = "2.357.13"
$ws.Range("E3").Value = "  -5.47%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'512.89"
$ws.Range("E5").Value = "  -3.90%  "
$ws.Range("D6").Value = "'127.44"
$ws.Range("E6").Value = "  -5.74%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.552"
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("D9").Value = "2.371.33"
$ws.Range("E9").Value = "  -5.85%  "
$ws.Range("D10").Value = "'0.0956"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("E12").Value = "  -8.78%  "
$ws.Range("E13").Value = "  -5.55%  "
$ws.Range("D14").Value = "2.772.01"
$ws.Range("E14").Value = "  -5.77%  "
$ws.Range("D15").Value = "56.484.31"
$ws.Range("E15").Value = "  -4.39%  "
$ws.Range("D16").Value = "'21.42"
$ws.Range("E16").Value = "  -4.66%  "
$ws.Range("E17").Value = "  -4.62%  "
$ws.Range("D18").Value = "2.331.13"
$ws.Range("E18").Value = "  -6.72%  "
$ws.Range("D19").Value = "'10.27"
$ws.Range("E19").Value = "  -4.04%  "
$ws.Range("D20").Value = "'4.04"
$ws.Range("E20").Value = "  -4.81%  "
$ws.Range("D21").Value = "'310.31"
$ws.Range("E21").Value = "  -3.62%  "
$ws.Range("D22").Value = "'6.09"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'65.13"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'0.388"
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").Value = "2.464.32"
$ws.Range("E27").Value = "  -5.94%  "
$ws.Range("E28").Value = "  -4.40%  "
$ws.Range("D29").Value = "'7.19"
$ws.Range("E29").Value = "  -4.27%  "
$ws.Range("D30").Value = "'174.82"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("D32").Value = "0.0₃0714"
$ws.Range("E32").Value = "  -7.21%  "
$ws.Range("D33").Value = "'6.12"
$ws.Range("E33").Value = "  -3.30%  "
$ws.Range("E34").Value = "  -7.05%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'17.62"
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("E38").Value = "  -5.98%  "
$ws.Range("D39").Value = "'3.71"
$ws.Range("E39").Value = "  -7.16%  "
$ws.Range("D40").Value = "'0.808"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("D41").Value = "'35.46"
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("E42").Value = "  -6.81%  "
$ws.Range("D43").Value = "'3.36"
$ws.Range("E43").Value = "  -3.78%  "
$ws.Range("D44").Value = "'4.87"
$ws.Range("E44").Value = "  -6.55%  "
$ws.Range("D45").Value = "'252.95"
$ws.Range("E45").Value = "  -9.90%  "
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0904"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'120.31"
$ws.Range("E48").Value = "  -8.78%  "
$ws.Range("D49").Value = "'0.0487"
$ws.Range("E49").Value = "  -4.43%  "
$ws.Range("D50").Value = "'0.0208"
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("D51").Value = "'16.57"
$ws.Range("E51").Value = "  -6.88%  "
